$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the cell values (English chemo/drug names) with the Chinese herb names,
# and extend the single column from 4 rows (with a gap at row 2) to 5 contiguous rows.
$ws.Range("A1").Value = "三七"
$ws.Range("A2").Value = "丝瓜络"
$ws.Range("A3").Value = "三棱"
$ws.Range("A4").Value = "丹参"
$ws.Range("A5").Value = "丹皮"

# A1 previously carried the underlined "hyperlink-ish" font style; drop it so the
# cell goes back to the default (unstyled) format.
$ws.Range("A1").ClearFormats()

# Flag duplicate herb names across the list with conditional formatting
# (Excel's built-in "Duplicate Values" rule / red text on light-red fill).
$ws.Range("A1:A5").FormatConditions.Delete()
$cf = $ws.Range("A1:A5").FormatConditions.AddUniqueValues()
$cf.DupeUnique = 1
$cf.Font.Color = 393372
$cf.Interior.Color = 13551615

# Move the active selection to B2.
[void]$ws.Range("B2").Select()
